# Generate Report for Handback
# Updates the handback-status workbook to reflect the new handoff/handback
# identifiers, xliff file names and timestamps produced by a later run.

$wb = $excel.ActiveWorkbook

$newGuid1 = "e7124884-0877-4844-a235-bea08737dd45"
$newGuid2 = "ffff65dcf7da-4ad3-4385-a7c7-8686eed17529"

$newHash1 = "1e2602b1917371dd72aa01bc3efb50038639f6c5"

# -------------------------------------------------------------------------
# Sheet 1: "Overview"
# -------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Cells.Item(2, 1).Value = "$newGuid1.md"
$wsOverview.Cells.Item(2, 2).Value = "e2e\$newGuid1.md"
$wsOverview.Cells.Item(2, 7).Value = "2016-09-07 01:22:41"

$wsOverview.Cells.Item(3, 1).Value = "$newGuid2.md"
$wsOverview.Cells.Item(3, 2).Value = "e2e\$newGuid2.md"
$wsOverview.Cells.Item(3, 7).Value = "2016-09-07 01:22:41"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($hl.Range.Row -eq 3) {
        $hl.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# -------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# -------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Cells.Item(2, 1).Value = "$newGuid1.md"
$wsZhCn.Cells.Item(2, 7).Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Cells.Item(2, 8).Value = "2016-09-07 01:22:36"
$wsZhCn.Cells.Item(2, 9).Value = "$newGuid1.md"
$wsZhCn.Cells.Item(2, 10).Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Cells.Item(2, 11).Value = "2016-09-07 01:22:53"

$wsZhCn.Cells.Item(3, 1).Value = "$newGuid2.md"
$wsZhCn.Cells.Item(3, 7).Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value = "2016-09-07 01:22:36"
$wsZhCn.Cells.Item(3, 9).Value = "$newGuid2.md"
$wsZhCn.Cells.Item(3, 10).Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 11).Value = "2016-09-07 01:22:53"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($hl.Range.Row -eq 3) {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}

# -------------------------------------------------------------------------
# Sheet 3: "de-de"
# -------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Cells.Item(2, 1).Value = "$newGuid1.md"
$wsDeDe.Cells.Item(2, 7).Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Cells.Item(2, 8).Value = "2016-09-07 01:22:41"
$wsDeDe.Cells.Item(2, 9).Value = "$newGuid1.md"
$wsDeDe.Cells.Item(2, 10).Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Cells.Item(2, 11).Value = "2016-09-07 01:23:02"

$wsDeDe.Cells.Item(3, 1).Value = "$newGuid2.md"
$wsDeDe.Cells.Item(3, 7).Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value = "2016-09-07 01:22:41"
$wsDeDe.Cells.Item(3, 9).Value = "$newGuid2.md"
$wsDeDe.Cells.Item(3, 10).Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDeDe.Cells.Item(3, 11).Value = "2016-09-07 01:23:02"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($hl.Range.Row -eq 3) {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}
